$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 3=0.04823258445971135; 4=0.1224246141273895; 5=0.1476342773254267; 6=2.059477635247433; 7=1.402221260560992; 8=1.264792942939124; 10=0.2025163487838313; 11=1.997077446537446; 14=1.337031879664575 }
    3 = @{ 3=0.04286911234592594; 4=0.1198512049305833; 5=0.1444842076393549; 6=2.045598682726435; 7=1.388397046314395; 8=1.265332143120119; 10=0.1980897476281811; 11=1.832599853824945; 14=1.354643802496216 }
    4 = @{ 3=0.03959357139387976; 4=0.1183083128177955; 5=0.1426247150590498; 6=2.038466716932092; 7=1.381007647272696; 8=1.266427981994497; 10=0.1954972664896601; 11=1.732298993553627; 14=1.366040991631554 }
    5 = @{ 3=0.03826311294965024; 4=0.1176889741765237; 5=0.1418857152039408; 6=2.035908665409295; 7=1.378271336674061; 8=1.267066374194343; 10=0.1944722517749895; 11=1.69159868694652; 14=1.370832048508017 }
    6 = @{ 3=0.03804245158517006; 4=0.117586702843937; 5=0.1417641373601057; 6=2.035504903862801; 7=1.377833540045316; 8=1.267183950856648; 10=0.1943039447830017; 11=1.68485088363127; 14=1.371636452225637 }
    7 = @{ 3=0.03957561086328099; 4=0.1182999220625263; 5=0.1426146727160429; 6=2.038430809551599; 7=1.380969633049148; 8=1.266435815510334; 10=0.1954833156144744; 11=1.731749394264227; 14=1.366105012008699 }
    8 = @{ 3=0.04637954083946738; 4=0.1215296125259471; 5=0.1465326256874953; 6=2.054402999223967; 7=1.397225805595724; 8=1.26481983891108; 10=0.2009639481352963; 11=1.940222402903942; 14=1.342983331999875 }
    9 = @{ 3=0.05986674108068257; 4=0.1281563473433636; 5=0.1548094828464741; 6=2.096810044606585; 7=1.437887059129991; 8=1.267742370350959; 10=0.212712738344365; 11=2.354535662301032; 14=1.302274490156506 }
    10 = @{ 3=0.06987143655841521; 4=0.1332020794362307; 5=0.1612552224944892; 6=2.134812701045632; 7=1.473211064882378; 8=1.273638032478885; 10=0.2219641849698348; 11=2.662361062089076; 14=1.275194536277638 }
    11 = @{ 3=0.07444527041840843; 4=0.1355356437042445; 5=0.1642674208564756; 6=2.153607338401272; 7=1.490485522160839; 8=1.277141886667749; 10=0.2263095358627112; 11=2.803162785800907; 14=1.263490554850737 }
    12 = @{ 3=0.07618063378271245; 4=0.1364247631849906; 5=0.1654196037325235; 6=2.160942579596735; 7=1.497201877956883; 8=1.278587498223544; 10=0.2279748292319965; 11=2.856592408569668; 14=1.259147115761518 }
    13 = @{ 3=0.07580674163897072; 4=0.1362330338061355; 5=0.1651709473820873; 6=2.159353080583102; 7=1.495747587110856; 8=1.278270867713701; 10=0.2276152956833783; 11=2.845080438748823; 14=1.260078609182521 }
    14 = @{ 3=0.07458797208448686; 4=0.1356086831629568; 5=0.1643619804349967; 6=2.154206433525658; 7=1.491034566673079; 8=1.277258434095273; 10=0.2264461429192437; 11=2.807556245347484; 14=1.26313144059214 }
    15 = @{ 3=0.07384187987598523; 4=0.1352269585940036; 5=0.1638679669633021; 6=2.151082411379235; 7=1.488170528458085; 8=1.276653775370931; 10=0.2257325859169299; 11=2.784586066465693; 14=1.265012932824533 }
    16 = @{ 3=0.06957298959052594; 4=0.133050340723571; 5=0.1610599801407062; 6=2.133614872865806; 7=1.472106527017502; 8=1.273425636230201; 10=0.221682970207226; 11=2.653174864131131; 14=1.275971807175878 }
    17 = @{ 3=0.0669600315939789; 4=0.1317248141604495; 5=0.1593578782501694; 6=2.123286089977341; 7=1.462561616973971; 8=1.271656200894796; 10=0.2192337961918582; 11=2.572755856434469; 14=1.282852327868614 }
    18 = @{ 3=0.06545925323088397; 4=0.1309660081137878; 5=0.1583864025898194; 6=2.117487011537534; 7=1.457185000002966; 8=1.270715792619399; 10=0.2178379665473216; 11=2.526573465214483; 14=1.286867679523283 }
    19 = @{ 3=0.06495147772207588; 4=0.1307097095343579; 5=0.158058769879311; 6=2.115547848684159; 7=1.455383992752843; 8=1.270410649128792; 10=0.2173675692461075; 11=2.51094934416011; 14=1.288237140037765 }
    20 = @{ 3=0.0672379647155168; 4=0.1318655463109053; 5=0.1595382905707439; 6=2.124370924367426; 7=1.463565945095723; 8=1.27183655332621; 10=0.21949318197602; 11=2.581309092200343; 14=1.282113895252294 }
    21 = @{ 3=0.07494586251851842; 4=0.1357919225739579; 5=0.1645992805459073; 6=2.155712197489962; 7=1.492414135945381; 8=1.277552582032968; 10=0.2267890129365497; 11=2.818574993819198; 14=1.262232343710558 }
    22 = @{ 3=0.08000300068181332; 4=0.1383897772062994; 5=0.1679741407406325; 6=2.177467567696311; 7=1.512288245814261; 8=1.281980946899381; 10=0.2316727501110307; 11=2.974289510642166; 14=1.249755152737784 }
    23 = @{ 3=0.07730208750825796; 4=0.137000364966525; 5=0.1661667560055804; 6=2.165739455458919; 7=1.501587196104708; 8=1.279553873758033; 10=0.229055597322926; 11=2.891122389988425; 14=1.256367133425805 }
    24 = @{ 3=0.06711230665626999; 4=0.1318019111586324; 5=0.1594567040935857; 6=2.123880037633569; 7=1.463111543092566; 8=1.271754776606969; 10=0.2193758755700799; 11=2.57744201371105; 14=1.28244755484695 }
    25 = @{ 3=0.05620183905715237; 4=0.1263324212145704; 5=0.1525065261315248; 6=2.084141974768031; 7=1.425936980850366; 8=1.26629603067363; 10=0.2094262408622569; 11=2.241856433415478; 14=1.312790906575486 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col]
    }
}
